$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Article 89 goes live: shift the rotating "ser" blog references in row 7.
# C7 gets the new article (89), E7 takes over C7's old value (88),
# I7 takes over E7's old value (87).
$ws.Range("C7").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 89"
$ws.Range("E7").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 88"
$ws.Range("I7").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 87"
